$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 15 entirely (mirrors the author selecting the whole row before
# moving it to the bottom of the list).
$ws.Range("A15:XFD15").Select()

# Remember the content + formatting of row 15 (the row being moved to the
# bottom of the participant list).
$movedA = $ws.Cells.Item(15, 1).Value2
$movedB = $ws.Cells.Item(15, 2).Value2
$ws.Cells.Item(15, 1).Copy()
$fmtRange = $ws.Range("A40:B40")
$fmtRange.PasteSpecial(-4122)
$ws.Cells.Item(15, 2).Copy()
$ws.Range("B40").PasteSpecial(-4122)

# Shift rows 16-25 up by one row (values and formatting together).
for ($r = 15; $r -le 24; $r++) {
    $src = $r + 1
    $valA = $ws.Cells.Item($src, 1).Value2
    $valB = $ws.Cells.Item($src, 2).Value2

    $ws.Cells.Item($src, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($src, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $valA
    $ws.Cells.Item($r, 2).Value = $valB
}

# Place the originally-row-15 content (and formatting) at the bottom, row 25.
$ws.Range("A40").Copy()
$ws.Cells.Item(25, 1).PasteSpecial(-4122)
$ws.Range("B40").Copy()
$ws.Cells.Item(25, 2).PasteSpecial(-4122)
$ws.Cells.Item(25, 1).Value = $movedA
$ws.Cells.Item(25, 2).Value = $movedB

# Clean up the scratch cells used to stash row 15's formatting.
$ws.Range("A40:B40").Clear()

# Restore the selection to row 15 (as recorded in the saved file).
$ws.Range("A15:XFD15").Select()
